$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1107593333333333
$ws.Range("N2").Value = 0.332278
$ws.Range("O2").Value = 0.1221364674089439
$ws.Range("P2").Value = 0.1221364674089439
$ws.Range("Q2").Value = 0.07810814642266667
$ws.Range("R2").Value = 0.7029733178040001
$ws.Range("S2").Value = 0.1221364674089439
$ws.Range("T2").Value = 0.1221364674089439

# Row 3 (Target cluster: FAPs)
$ws.Range("O3").Value = 0.4983148609452437
$ws.Range("P3").Value = 0.4983148609452437
$ws.Range("S3").Value = 0.4983148609452437
$ws.Range("T3").Value = 0.4983148609452437

# Row 4 (Target cluster: MuSCs)
$ws.Range("M4").Value = 0.3441933333333333
$ws.Range("N4").Value = 1.03258
$ws.Range("O4").Value = 0.3795486716458124
$ws.Range("P4").Value = 0.3795486716458124
$ws.Range("Q4").Value = 0.2427272038266666
$ws.Range("R4").Value = 2.18454483444
$ws.Range("S4").Value = 0.3795486716458124
$ws.Range("T4").Value = 0.3795486716458124
